$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "End" and "Target" variables were renamed everywhere:
#   Target -> Expected
#   End    -> Result
$ws.Range("B1").Value = "Expected"
$ws.Range("D1").Value = "Result"

# Update the active selection to match the saved workbook state.
$ws.Range("D10").Select() | Out-Null
